# The "Register" worksheet lists stakeholders and their Role (column E).
# Guozhi Yin (row 6) was "PM" and Cong Shang (row 7) was "Developer";
# both are unified under a single "Team member" role.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

$ws.Range("E6").Value = "Team member"
$ws.Range("E7").Value = "Team member"

